$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.406.28'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.579.39'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.46'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.490'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.47'
$ws.Range('E8').Value = '  -3.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.91'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0589'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '1.802.05'
$ws.Range('D14').Value = '1.576.22'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.517'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '28.414.66'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.76'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '231.91'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.07'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  +2.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.70'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0481'
$ws.Range('E31').Value = '  +3.88%  '
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').Value = '1.392.94'
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('E36').Value = '  +8.21%  '
$ws.Range('E37').Value = '  -2.93%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.65'
$ws.Range('E39').Value = '  +3.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.521'
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.90'
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.786'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('E45').Value = '  -3.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0455'
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.931'
$ws.Range('E47').Value = '  -5.05%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '62.58'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('D49').Value = '1.715.24'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.90'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '41.79'
$ws.Range('E51').Value = '  +5.32%  '
